$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 177; this shifts the existing rows
# (old 177-198) down to 179-200, matching the diff.
$ws.Rows("177:178").Insert()

# New row 177
$ws.Cells.Item(177, 1).Value = 6
$ws.Cells.Item(177, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(177, 3).Value = "Metropolitana"
$ws.Cells.Item(177, 4).Value = [DateTime]"2021-11-05"
$ws.Cells.Item(177, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = "Fruta"
$ws.Cells.Item(177, 7).Value = 100101
$ws.Cells.Item(177, 8).Value = "Berries"
$ws.Cells.Item(177, 9).Value = 100101001
$ws.Cells.Item(177, 10).Value = "Arándano (blue)"
$ws.Cells.Item(177, 11).Value = "Sin especificar"
$ws.Cells.Item(177, 12).Value = "Primera"
$ws.Cells.Item(177, 13).Value = 1350
$ws.Cells.Item(177, 14).Value = 7000
$ws.Cells.Item(177, 15).Value = 7000
$ws.Cells.Item(177, 16).Value = 7000
$ws.Cells.Item(177, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(177, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(177, 19).Value = 3500
$ws.Cells.Item(177, 20).Value = 2

# New row 178
$ws.Cells.Item(178, 1).Value = 6
$ws.Cells.Item(178, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(178, 3).Value = "Metropolitana"
$ws.Cells.Item(178, 4).Value = [DateTime]"2021-11-05"
$ws.Cells.Item(178, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(178, 5).Value = 13
$ws.Cells.Item(178, 6).Value = "Fruta"
$ws.Cells.Item(178, 7).Value = 100101
$ws.Cells.Item(178, 8).Value = "Berries"
$ws.Cells.Item(178, 9).Value = 100101001
$ws.Cells.Item(178, 10).Value = "Arándano (blue)"
$ws.Cells.Item(178, 11).Value = "Sin especificar"
$ws.Cells.Item(178, 12).Value = "Segunda"
$ws.Cells.Item(178, 13).Value = 100
$ws.Cells.Item(178, 14).Value = 6000
$ws.Cells.Item(178, 15).Value = 6000
$ws.Cells.Item(178, 16).Value = 6000
$ws.Cells.Item(178, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(178, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(178, 19).Value = 3000
$ws.Cells.Item(178, 20).Value = 2
